# Implement artificial intelligence for bot using minimax algorithm.
# This edit updates the UML class-diagram worksheet: renames/retypes several
# members, reorders the Board class method list, adds new BotAI methods,
# and adjusts the Game/Space classes accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1) Copy cell formatting (styles) onto the cells that need a *different*
#    style than they currently have, using stable "donor" cells whose
#    formatting we are not going to disturb until after we've copied from
#    them (A2/A1 never change; O5/A32/C20/C23 are copied from before their
#    own values/styles are rewritten below).
# ---------------------------------------------------------------------------

# Style "6" (italic, highlighted) donor -> targets that need style 6
$ws.Range("A32").Copy()
$ws.Range("A33").PasteSpecial($xlPasteFormats)
$ws.Range("A32").Copy()
$ws.Range("C15").PasteSpecial($xlPasteFormats)
$ws.Range("A32").Copy()
$ws.Range("C16").PasteSpecial($xlPasteFormats)
$ws.Range("A32").Copy()
$ws.Range("C21").PasteSpecial($xlPasteFormats)
$ws.Range("A32").Copy()
$ws.Range("C25").PasteSpecial($xlPasteFormats)
$ws.Range("A32").Copy()
$ws.Range("C28").PasteSpecial($xlPasteFormats)
$ws.Range("A32").Copy()
$ws.Range("E12").PasteSpecial($xlPasteFormats)

# Style "7" (new red highlighted font) for A15: start from style 6 (A32)
# then flip off italic and recolor the font red - matches the new font/xf
# that this revision introduces.
$ws.Range("A32").Copy()
$ws.Range("A15").PasteSpecial($xlPasteFormats)
$ws.Range("A15").Font.Italic = $false
$ws.Range("A15").Font.Color = 255

# Style "2" donor (O5, before it is rewritten) -> C20
$ws.Range("O5").Copy()
$ws.Range("C20").PasteSpecial($xlPasteFormats)

# Style "1" (bold header) donor -> targets that need style 1
$ws.Range("A1").Copy()
$ws.Range("C32").PasteSpecial($xlPasteFormats)
$ws.Range("A1").Copy()
$ws.Range("C33").PasteSpecial($xlPasteFormats)

# Style "4" (plain body) donor -> targets that need style 4
$ws.Range("A2").Copy()
$ws.Range("A32").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("A38").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("C23").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("C26").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("C27").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("C34").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("C35").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("C36").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("C37").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("O2").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("O3").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("O4").PasteSpecial($xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("O5").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Remove the two obsolete GameResult values that no longer fit after the
#    Board class method list grows (WinnerX/Winner0 at C30/C31 disappear;
#    their content moves up into C35/C36 below).
# ---------------------------------------------------------------------------
$ws.Range("C30").Clear()
$ws.Range("C31").Clear()

# ---------------------------------------------------------------------------
# 3) Write the final cell values (column A - Game class methods).
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = " +Game(board : Board)"
$ws.Range("A16").Value = " +Game()"
$ws.Range("A17").Value = " +GetBoard()"
$ws.Range("A18").Value = " +SetBoard(board : Board)"
$ws.Range("A19").Value = " +GetNewId()"
$ws.Range("A20").Value = " +GetCurrentTurnPlayer"
$ws.Range("A21").Value = " +SetCurrentTurnPlayer(Player)"
$ws.Range("A22").Value = " +GetNotCurrentTurnPlayer()"
$ws.Range("A23").Value = " +DetermineTurn()"
$ws.Range("A24").Value = " +GetUserPlayer()"
$ws.Range("A25").Value = " +SetUserPlayer(userShapeChoice)"
$ws.Range("A26").Value = " +GetBotPlayer"
$ws.Range("A27").Value = " +NewGame()"
$ws.Range("A28").Value = " +NewTurn()"
$ws.Range("A29").Value = " +UpdateScores(result : GameResult)"
$ws.Range("A30").Value = " +RestartGame()"
$ws.Range("A31").Value = " +GetUserShapeChoice()"
$ws.Range("A32").Value = " +PromptPickSpaceToOccupy() : Space"
$ws.Range("A33").Value = " +PromptPlayAgain() : bool"
$ws.Range("A34").Value = " +DisplayWinner(gameResult : GameResult)"
$ws.Range("A35").Value = " +DisplayPlayersScore()"
$ws.Range("A36").Value = " +GetPlayerFromShape(shape : Shape) : Player"
$ws.Range("A37").Value = " +GetPlayerFromResult(result : GameResult) : Player"
$ws.Range("A38").Value = " +ToString()"

# ---------------------------------------------------------------------------
# 4) Column C - Board class members/methods + Enum box.
# ---------------------------------------------------------------------------
$ws.Range("C15").Value = " +OccupySpace(board : Board, space Space, player : Player?)"
$ws.Range("C16").Value = " +OccupySpace(board : Board, space Space)"
$ws.Range("C17").Value = " +CheckWin() : (hasWinner : bool, winner : Shape?)"
$ws.Range("C18").Value = " +CheckTie() : bool"
$ws.Range("C19").Value = " +GetSpace(position) : Space"
$ws.Range("C20").Value = " +SetSpace(space : Space)"
$ws.Range("C21").Value = " +GetBoardClone(Board) : Board"
$ws.Range("C22").Value = " +PrintBoard()"
$ws.Range("C23").Value = " +GetBoardSpaceFromInt(int)"
$ws.Range("C24").Value = " +GetBoardSpaceFromCoordinates(x : int, y : int)"
$ws.Range("C25").Value = " +GetResultFromBoard(board)"
$ws.Range("C26").Value = " +GetResult()"
$ws.Range("C27").Value = " +SetResult()"
$ws.Range("C28").Value = " +GetShapeOfTurnFromBoard(Board board) : Shape"
$ws.Range("C29").Value = " +ToString()"
$ws.Range("C32").Value = "Enum"
$ws.Range("C33").Value = "BoardState"
$ws.Range("C34").Value = "Incomplete"
$ws.Range("C35").Value = "WinnerX"
$ws.Range("C36").Value = "Winner0"
$ws.Range("C37").Value = "Tie"

# ---------------------------------------------------------------------------
# 5) Column E - Space class updates.
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = " -occupant : Shape?"
$ws.Range("E7").Value = " +GetPosition() : Position"
$ws.Range("E8").Value = " +GetOccupant() : Player"
$ws.Range("E12").Value = " +GetSpaceClone(Space) : Space"

# ---------------------------------------------------------------------------
# 6) Column O - BotAI class (renamed from BotArtificialIntelligence) with
#    the new minimax-based methods.
# ---------------------------------------------------------------------------
$ws.Range("O1").Value = "BotAI"
$ws.Range("O2").Value = " +GetRandomMove(board Board) : Space"
$ws.Range("O3").Value = " +GetMinimaxMove(board : Board) : Space"
$ws.Range("O4").Value = " +Minimax(board : Board, isMaximizing : bool) : int, Space"
$ws.Range("O5").Value = " +GetScore(board : Board) : int"

# ---------------------------------------------------------------------------
# 7) Update the active selection to match the edited workbook (the author
#    ended up with the new BotAI method block selected).
# ---------------------------------------------------------------------------
$ws.Range("O2:O5").Select()
